{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the\n// document's table cells with its new value, per the commit's diff.\n// Uses Body.search(text, {matchCase:true}) to locate the exact literal\n// text run, then Range.insertText(newText, \"Replace\") to swap it in\n// place (keeps the run's existing formatting/rPr untouched).\nconst replacements = [\n  [\"83\u00d750=\", \"41\u00d745=\"],\n  [\"53\u00d787=\", \"23\u00d764=\"],\n  [\"29\u00d734=\", \"45\u00d758=\"],\n  [\"49\u00d726=\", \"81\u00d717=\"],\n  [\"57\u00d772=\", \"63\u00d758=\"],\n  [\"13\u00d730=\", \"47\u00d781=\"],\n  [\"80\u00d775=\", \"52\u00d796=\"],\n  [\"48\u00d771=\", \"39\u00d765=\"],\n  [\"37\u00d769=\", \"13\u00d732=\"],\n  [\"43\u00d775=\", \"68\u00d776=\"],\n  [\"16\u00d733=\", \"34\u00d759=\"],\n  [\"39\u00d715=\", \"42\u00d752=\"],\n  [\"29\u00d775=\", \"16\u00d794=\"],\n  [\"98\u00d786=\", \"37\u00d779=\"],\n  [\"69\u00d716=\", \"13\u00d798=\"],\n  [\"36\u00d767=\", \"25\u00d768=\"],\n  [\"32\u00d791=\", \"71\u00d754=\"],\n  [\"63\u00d734=\", \"80\u00d781=\"],\n  [\"23\u00d770=\", \"27\u00d759=\"],\n  [\"63\u00d720=\", \"58\u00d770=\"],\n  [\"80\u00d783=\", \"35\u00d714=\"],\n  [\"79\u00d792=\", \"36\u00d737=\"],\n  [\"96\u00d748=\", \"69\u00d739=\"],\n  [\"99\u00d739=\", \"96\u00d769=\"],\n  [\"47\u00d734=\", \"63\u00d756=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt in the\n# document's table cells with its new value, per the commit's diff.\n# Uses Find/Replace (wdReplaceOne semantics via Execute) scoped to the\n# whole document content; each old literal is unique in the document\n# so a single Execute per pair is safe and leaves formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"83\u00d750=\", \"41\u00d745=\"),\n    @(\"53\u00d787=\", \"23\u00d764=\"),\n    @(\"29\u00d734=\", \"45\u00d758=\"),\n    @(\"49\u00d726=\", \"81\u00d717=\"),\n    @(\"57\u00d772=\", \"63\u00d758=\"),\n    @(\"13\u00d730=\", \"47\u00d781=\"),\n    @(\"80\u00d775=\", \"52\u00d796=\"),\n    @(\"48\u00d771=\", \"39\u00d765=\"),\n    @(\"37\u00d769=\", \"13\u00d732=\"),\n    @(\"43\u00d775=\", \"68\u00d776=\"),\n    @(\"16\u00d733=\", \"34\u00d759=\"),\n    @(\"39\u00d715=\", \"42\u00d752=\"),\n    @(\"29\u00d775=\", \"16\u00d794=\"),\n    @(\"98\u00d786=\", \"37\u00d779=\"),\n    @(\"69\u00d716=\", \"13\u00d798=\"),\n    @(\"36\u00d767=\", \"25\u00d768=\"),\n    @(\"32\u00d791=\", \"71\u00d754=\"),\n    @(\"63\u00d734=\", \"80\u00d781=\"),\n    @(\"23\u00d770=\", \"27\u00d759=\"),\n    @(\"63\u00d720=\", \"58\u00d770=\"),\n    @(\"80\u00d783=\", \"35\u00d714=\"),\n    @(\"79\u00d792=\", \"36\u00d737=\"),\n    @(\"96\u00d748=\", \"69\u00d739=\"),\n    @(\"99\u00d739=\", \"96\u00d769=\"),\n    @(\"47\u00d734=\", \"63\u00d756=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceOne=1 (exactly one hit per old string)\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n"}
